# Generate Report for Handback
# Update the handback status report with refreshed timestamps/status for
# the 8fabe3d2-... (zh-cn) and f6cd9c42-... (de-de unaffected here) rows.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 00:18:03"
$wsOverview.Range("G5").Value = "2016-08-31 00:18:03"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Status column (E): "ht" -> "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H3").Value = "2016-08-31 00:17:56"
$wsZhCn.Range("H5").Value = "2016-08-31 00:17:56"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K3").Value = "2016-08-31 00:18:29"
$wsZhCn.Range("K5").Value = "2016-08-31 00:18:29"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Correspond Handoff Datetime column (H)
$wsDeDe.Range("H3").Value = "2016-08-31 00:18:03"
$wsDeDe.Range("H5").Value = "2016-08-31 00:18:03"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K3").Value = "2016-08-31 00:18:36"
$wsDeDe.Range("K5").Value = "2016-08-31 00:18:36"
